# Update the Answer column (C) values per the sourced query results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: "ANKITA" -> "Reliance Retail Limited."
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = 'Reliance Retail Limited.'

# C3: "ANKITA" -> "iPhone XR 128GB"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = 'iPhone XR 128GB'

# C4: "Google" -> "iPhone"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = 'iPhone'

# C5: "CNB 30W Fast Charging Type C Rapidly Ada..." -> "iPhone XR 128GB..."
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'iPhone XR 128GB'

# C6: "2025-01-05" -> "2019-08-03"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '2019-08-03'

# C7: "05.01.2025" -> "2019-08-03"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = '2019-08-03'

# C8: "07AZAPA3803E1Z5" -> "29AABCR1718E1ZL"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = '29AABCR1718E1ZL'

# C9: "IN-2529" -> "8884136002703082019"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = '8884136002703082019'

# C11: "₹699.00" -> "The total amount is Rs 63591.22."
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'The total amount is Rs 63591.22.'

# C12: "No discount provided." -> "Discount: Rs 1297.78"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'Discount: Rs 1297.78'

# C13: "₹106.63" -> "15.08"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '15.08'

# C14: "18%" -> "0.00%"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '0.00%'

# C15: "ANKITA* 4649/C -75, Street No. 10, New M..." -> "#LG10, Phoenix Market City, Opp. Mahadev..."
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '#LG10, Phoenix Market City, Opp. Mahadevapura CMC Office, Whitefield Road, Bangalore - 560048.'

# C16: "1/57/A/194 , sri Aurobindo Residency, 9,..." -> "Customer Address:V SRINIVASUFlat No: 404..."
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'Customer Address:V SRINIVASUFlat No: 404Flr No: 4TH FLOORWing: B BLOCKBldg: ALEMBIC URBAN FORESTSoc: AB VAJPAYEE ROADSec/Loc: .Plot No: CHENNASANDRA MAIN ROADStreet: OPP TO WHITEFIELD GLOBAL SCHOOLArea: KADUGODICity: BANGALORE EASTState: KARNATAKAPinCode: 560067Contact# 9886642984NEELIMA_VE@YAHOO.CO.INRelationship ID: 9886642984'

# C17: "LoukyaFlat 302 , mahadev enclave, Janaki..." -> "Customer Address:V SRINIVASUFlat No: 404..."
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'Customer Address:V SRINIVASUFlat No: 404Flr No: 4TH FLOORWing: B BLOCKBldg: ALEMBIC URBAN FORESTSoc: AB VAJPAYEE ROADSec/Loc: .Plot No: CHENNASANDRA MAIN ROADStreet: OPP TO WHITEFIELD GLOBAL SCHOOLArea: KADUGODICity: BANGALORE EASTState: KARNATAKAPinCode: 560067'

# C18: "ANKITA* 4649/C -75, Street No. 10, New M..." -> "#LG10, Phoenix Market City, Opp. Mahadev..."
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = '#LG10, Phoenix Market City, Opp. Mahadevapura CMC Office, Whitefield Road, Bangalore - 560048.'

# C19: "Syam" -> "V SRINIVASU"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'V SRINIVASU'

# C20: "AZAPA3803E" -> "U01100MH1999PLC120563"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'U01100MH1999PLC120563'

# C21: "UPI" -> "EMI"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'EMI'

# C22: "CNB 30W Fast Charging Type C Rapidly Ada..." -> "iPhone XR 128GB..."
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'iPhone XR 128GB'

# C23: "2026-01-05" -> "2020-08-02"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = '2020-08-02'
